$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell-value updates derived from the authoritative diff (cryptos.xlsx,
# GitHub Actions refresh). Each assignment mirrors one inline-string <t>
# change in the OOXML diff.
#
# D-column values such as "347.74" or "0.550" are numeric-looking text
# (thousand-dot formatted prices / fixed decimal-precision quotes) that
# must stay literal text -- a plain .Value assignment would let Excel
# parse them into floating point numbers and silently drop meaningful
# trailing zeros (e.g. "0.550" -> 0.55). For those cells we briefly flip
# NumberFormat to Text ("@") before writing, then restore the original
# "Normal" style so no stray formatting change is left behind.

$ws.Range("D2").Value = "52.289.31"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "2.794.21"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "115.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.550"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.79%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D15").Value = "3.232.06"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "2.772.49"
$ws.Range("E16").Value = "  +1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.890"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "52.214.79"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  +7.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.19%  "
$ws.Range("E21").Value = "  -3.41%  "
$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.22%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  +3.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.79"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "50.21"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0448"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +28.98%  "
$ws.Range("E34").Value = "  +2.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0831"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.02%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.59"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.03%  "
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("E41").Value = "  +9.36%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.115"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.94%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").Value = "2.056.34"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E48").Value = "  +3.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.966"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.08%  "
